$wb = $excel.ActiveWorkbook

# Rename sheets (by current index/position, since names change)
$wb.Worksheets.Item(1).Name = "GNG_TO-16509962017544417"
$wb.Worksheets.Item(2).Name = "NB_TO-16509962039257438"
$wb.Worksheets.Item(3).Name = "RS_TO-16509962039257438"
$wb.Worksheets.Item(4).Name = "TOL_TO-1650996203973741"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509962040377407"

# Sheet 1: GNG_TO
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509962017224195.csv"
$ws1.Range("B3").Value = "GNG_stims-16509962017384055.csv"
$ws1.Range("B4").Value = "go_stims-16509962017384055.csv"
$ws1.Range("B5").Value = "GNG_stims-16509962017544417.csv"

# Sheet 2: NB_TO
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16509962030697036.csv"
$ws2.Range("B3").Value = "ZB-match_6-16509962026455684.csv"
$ws2.Range("B4").Value = "TB-16509962038057024.csv"
$ws2.Range("B5").Value = "OB-16509962027096055.csv"
$ws2.Range("B6").Value = "TB-16509962035576982.csv"
$ws2.Range("B7").Value = "TB-16509962039097009.csv"
$ws2.Range("B8").Value = "OB-16509962034297295.csv"
$ws2.Range("B9").Value = "ZB-match_4-16509962019736094.csv"
$ws2.Range("B10").Value = "ZB-match_6-16509962025175674.csv"

# Sheet 3: RS_TO - no changes

# Sheet 4: TOL_TO
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1650996203941738.csv"
$ws4.Range("B3").Value = "ZM_stims-16509962039257438.csv"
$ws4.Range("B4").Value = "MM_stims-16509962039577115.csv"
$ws4.Range("B5").Value = "ZM_stims-1650996203941738.csv"
$ws4.Range("B6").Value = "MM_stims-1650996203973741.csv"
$ws4.Range("B7").Value = "ZM_stims-16509962039577115.csv"

# Sheet 5: vSAT_TO
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1650996203973741.csv"
$ws5.Range("B3").Value = "SAT_stims-16509962039897368.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509962040057366.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509962040217361.csv"
